$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update name and email (remain text)
$ws.Range("A5").Value = "Victor Nascimento"
$ws.Range("E5").Value = "victornascimento@gmail.com"

# B5, C5, D5 switch from numeric to text values - force text storage
$ws.Range("B5:D5").NumberFormat = "@"
$ws.Range("B5").Value = "934856"
$ws.Range("C5").Value = "9487"
$ws.Range("D5").Value = "4987"

# F5 stays numeric, just update its value
$ws.Range("F5").Value = 344557543
